$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy formatting for the new G column (header + data) from column F,
#     before we touch anything else, so the new cells inherit the correct
#     cell styles (bold/border header style, wrap-text data style).
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("F2").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Remove the old data rows 3-5, keep only the header row and a single
#     data row (row 2), shifting the remaining rows up.
$ws.Range("A3:A5").EntireRow.Delete()

# --- Update the header row (row 1) ---
$ws.Range("A1").Value = "Subarea"
$ws.Range("B1").Value = "Location"
$ws.Range("C1").Value = "Approved CCTV Vendor"
$ws.Range("D1").Value = "Walkthrough"
$ws.Range("E1").Value = "Vendor Surveyor"
$ws.Range("F1").Value = "Reviewer"
$ws.Range("G1").Value = "Notes"

# --- Update the data row (row 2) ---
$ws.Range("A2").Value = "gh"
$ws.Range("B2").Value = "gh"
$ws.Range("C2").Value = "GPH"
$ws.Range("D2").Value = "Accepted"
$ws.Range("E2").Value = "gh"
$ws.Range("F2").Value = "T. Martin"
$ws.Range("G2").Value = "gh"

# --- Column widths: extend the 12.6640625-wide formatting through column H,
#     and give column I a wider 25.6640625 width for the Notes/legend area.
$ws.Columns.Item(7).ColumnWidth = 11.83
$ws.Columns.Item(8).ColumnWidth = 11.83
$ws.Columns.Item(9).ColumnWidth = 24.83

# --- Move the conditional "not blank" formatting from B2:I5 to B2:L2 ---
$fcs = $ws.Range("B2").FormatConditions
for ($i = 1; $i -le $fcs.Count(); $i++) {
    $fc = $fcs.Item($i)
    if ($fc.Type() -eq 10) {
        $fc.ModifyAppliesToRange($ws.Range("B2:L2"))
    }
}
